$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "300.78"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.16%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.45"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.60%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.157"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.75%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07358"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.87%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.819"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "23.15%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.832"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.27%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.755"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.42%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9307"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.19%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1699"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.12%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07031"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-9.21%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08144"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.61%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.75%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09945"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.51%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001490"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.09%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006154"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-4.31%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.458"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.74%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.221"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.30%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.80%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1331"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.53%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.560"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.16%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04642"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.81%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.1582"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-2.42%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001216"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.15%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004749"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "7.50%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001298"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-7.16%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "7.48%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01718"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-1.59%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04515"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.38%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007112"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.15%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1341"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.11%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002196"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.36%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01047"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-17.57%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006237"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.46%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-21.45%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.7398"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-60.49%"
